$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '69.687.36'
$ws.Range("E2").Value = '  +6.17%  '

$ws.Range("D3").Value = '3.587.18'
$ws.Range("E3").Value = '  +5.67%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.644'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.07%  '

$ws.Range("D8").Value = '3.580.78'
$ws.Range("E8").Value = '  +5.68%  '

$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.182'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.663'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000292'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.96%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.63%  '

$ws.Range("D15").Value = '4.155.68'
$ws.Range("E15").Value = '  +5.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.94%  '

$ws.Range("D17").Value = '3.585.00'
$ws.Range("E17").Value = '  +5.23%  '

$ws.Range("D18").Value = '69.658.72'
$ws.Range("E18").Value = '  +6.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.21%  '

$ws.Range("E20").Value = '  +1.12%  '

$ws.Range("E21").Value = '  +5.42%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '500.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +18.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.18'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.94%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.50'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +14.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '614.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '65.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.116'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.16%  '

$ws.Range("D36").Value = '0.0₃0835'
$ws.Range("E36").Value = '  +11.69%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.148'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.31%  '

$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.04%  '

$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.01'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.91%  '

$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.398'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.82%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.96%  '

$ws.Range("D42").Value = '3.341.33'
$ws.Range("E42").Value = '  +8.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.10'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.97%  '

$ws.Range("E44").Value = '  +10.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0444'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +16.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.99%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.138'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.24'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.71%  '

$ws.Range("E51").Value = '  +0.21%  '
